$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 14-16 (old rows corresponding to the removed "Resolving-Mac" target cluster)
$ws.Rows.Item(14).Resize(3).Delete() | Out-Null

$data = @(
    @("ECs", "Angpt1", "Tek", "ECs", 1, 0.3333333333333333, 0.2556606666666667, 0.766982, 0.01354513404628681, 0.01354513404628681, 3, 1, 71.47459166666665, 214.423775, 0.9365945886638486, 0.9365945886638485, 18.27324175522778, 164.45917579705, 0.01268629925047868, 0.01268629925047868),
    @("ECs", "Angpt1", "Tek", "FAPs", 1, 0.3333333333333333, 0.2556606666666667, 0.766982, 0.01354513404628681, 0.01354513404628681, 3, 1, 3.258629666666666, 9.775889, 0.04270069742396077, 0.04270069742396076, 0.8331034329997778, 7.497930896998, 0.0005783866704774823, 0.0005783866704774822),
    @("ECs", "Angpt1", "Tek", "Inflammatory-Mac", 1, 0.3333333333333333, 0.2556606666666667, 0.766982, 0.01354513404628681, 0.01354513404628681, 1, 0.3333333333333333, 0.1062743333333333, 0.318823, 0.001392606284175224, 0.001392606284175224, 0.02717016690955556, 0.244531502186, 0.00001886303879285478, 0.00001886303879285478),
    @("ECs", "Angpt1", "Tek", "MuSCs", 1, 0.3333333333333333, 0.2556606666666667, 0.766982, 0.01354513404628681, 0.01354513404628681, 3, 1, 1.47377, 4.42131, 0.01931210762801542, 0.01931210762801541, 0.3767850207133334, 3.39106518642, 0.0002615850865377868, 0.0002615850865377867),
    @("FAPs", "Angpt1", "Tek", "ECs", 3, 1, 16.677026, 50.03107799999999, 0.8835639662863414, 0.8835639662863415, 3, 1, 71.47459166666665, 214.423775, 0.9365945886638486, 0.9365945886638485, 1191.983623564383, 10727.85261207945, 0.8275412295621545, 0.8275412295621545),
    @("FAPs", "Angpt1", "Tek", "FAPs", 3, 1, 16.677026, 50.03107799999999, 0.8835639662863414, 0.8835639662863415, 3, 1, 3.258629666666666, 9.775889, 0.04270069742396077, 0.04270069742396076, 54.34425167537132, 489.0982650783419, 0.03772879757910774, 0.03772879757910774),
    @("FAPs", "Angpt1", "Tek", "Inflammatory-Mac", 3, 1, 16.677026, 50.03107799999999, 0.8835639662863414, 0.8835639662863415, 1, 0.3333333333333333, 0.1062743333333333, 0.318823, 0.001392606284175224, 0.001392606284175224, 1.772339820132667, 15.951058381194, 0.001230456731921145, 0.001230456731921145),
    @("FAPs", "Angpt1", "Tek", "MuSCs", 3, 1, 16.677026, 50.03107799999999, 0.8835639662863414, 0.8835639662863415, 3, 1, 1.47377, 4.42131, 0.01931210762801542, 0.01931210762801541, 24.57810060802, 221.20290547218, 0.01706348241315801, 0.01706348241315801),
    @("MuSCs", "Angpt1", "Tek", "ECs", 3, 1, 1.942037333333333, 5.826112, 0.1028908996673717, 0.1028908996673717, 3, 1, 71.47459166666665, 214.423775, 0.9365945886638486, 0.9365945886638485, 138.8063254014222, 1249.2569286128, 0.09636705985121535, 0.09636705985121535),
    @("MuSCs", "Angpt1", "Tek", "FAPs", 3, 1, 1.942037333333333, 5.826112, 0.1028908996673717, 0.1028908996673717, 3, 1, 3.258629666666666, 9.775889, 0.04270069742396077, 0.04270069742396076, 6.328380468174222, 56.955424213568, 0.004393513174375546, 0.004393513174375546),
    @("MuSCs", "Angpt1", "Tek", "Inflammatory-Mac", 3, 1, 1.942037333333333, 5.826112, 0.1028908996673717, 0.1028908996673717, 1, 0.3333333333333333, 0.1062743333333333, 0.318823, 0.001392606284175224, 0.001392606284175224, 0.2063887229084445, 1.857498506176, 0.0001432865134612243, 0.0001432865134612243),
    @("MuSCs", "Angpt1", "Tek", "MuSCs", 3, 1, 1.942037333333333, 5.826112, 0.1028908996673717, 0.1028908996673717, 3, 1, 1.47377, 4.42131, 0.01931210762801542, 0.01931210762801541, 2.862116360746667, 25.75904724672, 0.001987040128319619, 0.001987040128319619),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
